$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Forms export sync corrected the Q2_1/Q2_2/Q2_3 answers recorded for
# response row 11 (sheet row 11) from "Excellent 5" to " Bad 1"
# (the source value begins with a non-breaking space, matching the other
# rows that already use that rating).
$badOne = [string]([char]0x00A0) + "Bad 1"

$ws.Range("F11").Value = $badOne
$ws.Range("G11").Value = $badOne
$ws.Range("H11").Value = $badOne

# Reflect the cell that was last touched/selected when the file was saved.
[void]$ws.Range("H11").Select()
